# Case_1_248 res_bus/vm_pu.xlsx update: rerun with 380 kV slack setpoint (1.02 pu).
# Updates bus voltage magnitudes in columns B-F and I-N for rows 2-25 (row 1 is the
# header and columns G/H are unaffected).
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Row 2
$ws.Range("B2").Value = 1.02
$ws.Range("C2").Value = 1.024609252641882
$ws.Range("D2").Value = 1.031550077435396
$ws.Range("E2").Value = 1.025069889366229
$ws.Range("F2").Value = 1.035764980930472
$ws.Range("I2").Value = 1.033816135422419
$ws.Range("J2").Value = 1.029782908912855
$ws.Range("K2").Value = 1.034357871925654
$ws.Range("L2").Value = 1.027896538741618
$ws.Range("M2").Value = 1.038560649235037
$ws.Range("N2").Value = 1.031245318761266

# Row 3
$ws.Range("B3").Value = 1.02
$ws.Range("C3").Value = 1.025519766115473
$ws.Range("D3").Value = 1.032028970098973
$ws.Range("E3").Value = 1.02584193080756
$ws.Range("F3").Value = 1.036905067889644
$ws.Range("I3").Value = 1.033986102177138
$ws.Range("J3").Value = 1.03033260366658
$ws.Range("K3").Value = 1.034646235768816
$ws.Range("L3").Value = 1.02847590082571
$ws.Range("M3").Value = 1.039509320476909
$ws.Range("N3").Value = 1.031795794144588

# Row 4
$ws.Range("B4").Value = 1.02
$ws.Range("C4").Value = 1.026109356826506
$ws.Range("D4").Value = 1.032339136029018
$ws.Range("E4").Value = 1.026342241533654
$ws.Range("F4").Value = 1.037643502423261
$ws.Range("I4").Value = 1.034095068309438
$ws.Range("J4").Value = 1.030688110142315
$ws.Range("K4").Value = 1.034832394066635
$ws.Range("L4").Value = 1.028850877591822
$ws.Range("M4").Value = 1.040123317964751
$ws.Range("N4").Value = 1.032151805480296

# Row 5
$ws.Range("B5").Value = 1.02
$ws.Range("C5").Value = 1.026357322231704
$ws.Range("D5").Value = 1.032469597673264
$ws.Range("E5").Value = 1.026552750024266
$ws.Range("F5").Value = 1.037954112843225
$ws.Range("I5").Value = 1.034140634453223
$ws.Range("J5").Value = 1.030837520483723
$ws.Range("K5").Value = 1.034910550513087
$ws.Range("L5").Value = 1.029008538616078
$ws.Range("M5").Value = 1.040381476516981
$ws.Range("N5").Value = 1.032301428001526

# Row 6
$ws.Range("B6").Value = 1.02
$ws.Range("C6").Value = 1.026398962638525
$ws.Range("D6").Value = 1.032491506690992
$ws.Range("E6").Value = 1.026588105718766
$ws.Range("F6").Value = 1.038006275808441
$ws.Range("I6").Value = 1.034148270931422
$ws.Range("J6").Value = 1.030862604499377
$ws.Range("K6").Value = 1.034923667167575
$ws.Range("L6").Value = 1.029035011787843
$ws.Range("M6").Value = 1.040424824456391
$ws.Range("N6").Value = 1.03232654763936

# Row 7
$ws.Range("B7").Value = 1.02
$ws.Range("C7").Value = 1.026112669754795
$ws.Range("D7").Value = 1.032340878999198
$ws.Range("E7").Value = 1.026345053661143
$ws.Range("F7").Value = 1.037647652137305
$ws.Range("I7").Value = 1.034095678122824
$ws.Range("J7").Value = 1.03069010674681
$ws.Range("K7").Value = 1.034833438808763
$ws.Range("L7").Value = 1.02885298418562
$ws.Range("M7").Value = 1.040126767358834
$ws.Range("N7").Value = 1.032153804920198

# Row 8
$ws.Range("B8").Value = 1.02
$ws.Range("C8").Value = 1.02491687585162
$ws.Range("D8").Value = 1.031711860358719
$ws.Range("E8").Value = 1.02533064867825
$ws.Range("F8").Value = 1.036150128962936
$ws.Range("I8").Value = 1.033873786061641
$ws.Range("J8").Value = 1.029968718181677
$ws.Range("K8").Value = 1.034455414773418
$ws.Range("L8").Value = 1.028092317516147
$ws.Range("M8").Value = 1.038881227110102
$ws.Range("N8").Value = 1.031431391900558

# Row 9
$ws.Range("B9").Value = 1.02
$ws.Range("C9").Value = 1.022813048909476
$ws.Range("D9").Value = 1.030605744510038
$ws.Range("E9").Value = 1.023548924971929
$ws.Range("F9").Value = 1.033516845456389
$ws.Range("I9").Value = 1.033475043440816
$ws.Range("J9").Value = 1.028696175647258
$ws.Range("K9").Value = 1.033786020040151
$ws.Range("L9").Value = 1.026752662002594
$ws.Range("M9").Value = 1.036687544938164
$ws.Range("N9").Value = 1.030157042209778

# Row 10
$ws.Range("B10").Value = 1.02
$ws.Range("C10").Value = 1.021412778450242
$ws.Range("D10").Value = 1.029869973806245
$ws.Range("E10").Value = 1.022365075831965
$ws.Range("F10").Value = 1.031765066423194
$ws.Range("I10").Value = 1.033204042143105
$ws.Range("J10").Value = 1.027846945132532
$ws.Range("K10").Value = 1.033337618171422
$ws.Range("L10").Value = 1.02586010835458
$ws.Range("M10").Value = 1.035225864468192
$ws.Range("N10").Value = 1.029306605690313

# Row 11
$ws.Range("B11").Value = 1.02
$ws.Range("C11").Value = 1.020806997438605
$ws.Range("D11").Value = 1.029551786326544
$ws.Range("E11").Value = 1.021853412636592
$ws.Range("F11").Value = 1.031007417202179
$ws.Range("I11").Value = 1.03308547516639
$ws.Range("J11").Value = 1.027479022749476
$ws.Range("K11").Value = 1.033142960252397
$ws.Range("L11").Value = 1.025473764599098
$ws.Range("M11").Value = 1.034593128309253
$ws.Range("N11").Value = 1.028938160815272

# Row 12
$ws.Range("B12").Value = 1.02
$ws.Range("C12").Value = 1.02058206612613
$ws.Range("D12").Value = 1.029433659994602
$ws.Range("E12").Value = 1.021663502312045
$ws.Range("F12").Value = 1.030726125266673
$ws.Range("I12").Value = 1.033041251076737
$ws.Range("J12").Value = 1.027342330712366
$ws.Range("K12").Value = 1.033070581956895
$ws.Range("L12").Value = 1.025330280927125
$ws.Range("M12").Value = 1.034358129371736
$ws.Range("N12").Value = 1.028801274659791

# Row 13
$ws.Range("B13").Value = 1.02
$ws.Range("C13").Value = 1.020630310899767
$ws.Range("D13").Value = 1.029458995631558
$ws.Range("E13").Value = 1.021704232177412
$ws.Range("F13").Value = 1.030786457334039
$ws.Range("I13").Value = 1.03305074556878
$ws.Range("J13").Value = 1.027371652931206
$ws.Range("K13").Value = 1.033086110665253
$ws.Range("L13").Value = 1.025361057670088
$ws.Range("M13").Value = 1.034408536195855
$ws.Range("N13").Value = 1.028830638519544

# Row 14
$ws.Range("B14").Value = 1.02
$ws.Range("C14").Value = 1.020788402853904
$ws.Range("D14").Value = 1.029542020679198
$ws.Range("E14").Value = 1.021837711641731
$ws.Range("F14").Value = 1.030984162814942
$ws.Range("I14").Value = 1.033081823320311
$ws.Range("J14").Value = 1.027467724335996
$ws.Range("K14").Value = 1.033136978937225
$ws.Range("L14").Value = 1.025461903748226
$ws.Range("M14").Value = 1.034573702662548
$ws.Range("N14").Value = 1.028926846356748

# Row 15
$ws.Range("B15").Value = 1.02
$ws.Range("C15").Value = 1.02088581951067
$ws.Range("D15").Value = 1.02959318350436
$ws.Range("E15").Value = 1.021919971886016
$ws.Range("F15").Value = 1.031105993276859
$ws.Range("I15").Value = 1.033100947109752
$ws.Range("J15").Value = 1.027526913241818
$ws.Range("K15").Value = 1.033168310824889
$ws.Range("L15").Value = 1.025524041234191
$ws.Range("M15").Value = 1.034675470817556
$ws.Range("N15").Value = 1.028986119317606

# Row 16
$ws.Range("B16").Value = 1.02
$ws.Range("C16").Value = 1.021452993698449
$ws.Range("D16").Value = 1.029891099558619
$ws.Range("E16").Value = 1.022399053382245
$ws.Range("F16").Value = 1.03181536771259
$ws.Range("I16").Value = 1.03321188534727
$ws.Range("J16").Value = 1.027871358781194
$ws.Range("K16").Value = 1.03335052659029
$ws.Range("L16").Value = 1.025885751697692
$ws.Range("M16").Value = 1.035267860919046
$ws.Range("N16").Value = 1.029331054009156

# Row 17
$ws.Range("B17").Value = 1.02
$ws.Range("C17").Value = 1.021808913645538
$ws.Range("D17").Value = 1.030078084320943
$ws.Range("E17").Value = 1.022699824116686
$ws.Range("F17").Value = 1.032260575943535
$ws.Range("I17").Value = 1.033281147261461
$ws.Range("J17").Value = 1.028087367367801
$ws.Range("K17").Value = 1.033464693373511
$ws.Range("L17").Value = 1.026112680652366
$ws.Range("M17").Value = 1.03563950023462
$ws.Range("N17").Value = 1.029547369352733

# Row 18
$ws.Range("B18").Value = 1.02
$ws.Range("C18").Value = 1.022016568292506
$ws.Range("D18").Value = 1.03018718844051
$ws.Range("E18").Value = 1.022875350300015
$ws.Range("F18").Value = 1.032520343528521
$ws.Range("I18").Value = 1.033321428661346
$ws.Range("J18").Value = 1.028213342176921
$ws.Range("K18").Value = 1.033531236912221
$ws.Range("L18").Value = 1.026245057714607
$ws.Range("M18").Value = 1.035856288620335
$ws.Range("N18").Value = 1.029673523060531

# Row 19
$ws.Range("B19").Value = 1.02
$ws.Range("C19").Value = 1.022087382052343
$ws.Range("D19").Value = 1.030224396702491
$ws.Range("E19").Value = 1.022935215708925
$ws.Range("F19").Value = 1.032608932002731
$ws.Range("I19").Value = 1.033335143564169
$ws.Range("J19").Value = 1.028256293027005
$ws.Range("K19").Value = 1.03355391836523
$ws.Range("L19").Value = 1.026290197091003
$ws.Range("M19").Value = 1.035930210812202
$ws.Range("N19").Value = 1.029716534905749

# Row 20
$ws.Range("B20").Value = 1.02
$ws.Range("C20").Value = 1.02177072136092
$ws.Range("D20").Value = 1.030058018572279
$ws.Range("E20").Value = 1.022667544769094
$ws.Range("F20").Value = 1.032212800504768
$ws.Range("I20").Value = 1.033273728305182
$ws.Range("J20").Value = 1.028064193686156
$ws.Range("K20").Value = 1.033452449317379
$ws.Range("L20").Value = 1.026088331947011
$ws.Range("M20").Value = 1.035599625045433
$ws.Range("N20").Value = 1.029524162761801

# Row 21
$ws.Range("B21").Value = 1.02
$ws.Range("C21").Value = 1.02074184642364
$ws.Range("D21").Value = 1.029517570128889
$ws.Range("E21").Value = 1.021798401276625
$ws.Range("F21").Value = 1.030925939813872
$ws.Range("I21").Value = 1.033072676745725
$ws.Range("J21").Value = 1.027439434504025
$ws.Range("K21").Value = 1.033122001524951
$ws.Range("L21").Value = 1.025432206488315
$ws.Range("M21").Value = 1.034525064495281
$ws.Range("N21").Value = 1.028898516349971

# Row 22
$ws.Range("B22").Value = 1.02
$ws.Range("C22").Value = 1.020095430948906
$ws.Range("D22").Value = 1.029178132118423
$ws.Range("E22").Value = 1.021252770994926
$ws.Range("F22").Value = 1.030117607355595
$ws.Range("I22").Value = 1.032945208913142
$ws.Range("J22").Value = 1.027046454568418
$ws.Range("K22").Value = 1.032913810093766
$ws.Range("L22").Value = 1.025019799786102
$ws.Range("M22").Value = 1.033849604532536
$ws.Range("N22").Value = 1.02850497833778

# Row 23
$ws.Range("B23").Value = 1.02
$ws.Range("C23").Value = 1.02043806234855
$ws.Range("D23").Value = 1.029358039635562
$ws.Range("E23").Value = 1.021541940392261
$ws.Range("F23").Value = 1.030546046856303
$ws.Range("I23").Value = 1.033012882211727
$ws.Range("J23").Value = 1.027254796436706
$ws.Range("K23").Value = 1.033024216342542
$ws.Range("L23").Value = 1.025238412197397
$ws.Range("M23").Value = 1.034207663497773
$ws.Range("N23").Value = 1.028713616075419

# Row 24
$ws.Range("B24").Value = 1.02
$ws.Range("C24").Value = 1.021787978647993
$ws.Range("D24").Value = 1.030067085299869
$ws.Range("E24").Value = 1.022682130133603
$ws.Range("F24").Value = 1.032234387905838
$ws.Range("I24").Value = 1.033277080976526
$ws.Range("J24").Value = 1.028074664935219
$ws.Range("K24").Value = 1.03345798202778
$ws.Range("L24").Value = 1.026099334038044
$ws.Range("M24").Value = 1.035617642872981
$ws.Range("N24").Value = 1.029534648881239

# Row 25
$ws.Range("B25").Value = 1.02
$ws.Range("C25").Value = 1.02335654042687
$ws.Range("D25").Value = 1.030891419981353
$ws.Range("E25").Value = 1.024008849629563
$ws.Range("F25").Value = 1.034196953777262
$ws.Range("I25").Value = 1.033579042022276
$ws.Range("J25").Value = 1.029025315053099
$ws.Range("K25").Value = 1.033959455865145
$ws.Range("L25").Value = 1.027098902141314
$ws.Range("M25").Value = 1.037254529516405
$ws.Range("N25").Value = 1.030486649031328
